# Rename "site" to "option" across the workbook (header + site_1..site_4 values).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # LookAt:=xlPart (2) so "site_1".."site_4" (which contain "site" as a
    # substring) are updated along with the standalone "site" header.
    $ws.Cells.Replace("site", "option", 2, 1, $false, $false, $false, $false)
}
